$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.916.95'
$ws.Range("E2").Value = '  +1.09%  '

$ws.Range("D3").Value = '2.219.83'
$ws.Range("E3").Value = '  +1.00%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''250.82'
$ws.Range("E5").Value = '  -0.90%  '

$ws.Range("D6").Value = '''0.623'
$ws.Range("E6").Value = '  -0.57%  '

$ws.Range("D7").Value = '''67.93'
$ws.Range("E7").Value = '  -1.78%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '''0.633'
$ws.Range("E9").Value = '  +7.62%  '

$ws.Range("D10").Value = '''39.34'
$ws.Range("E10").Value = '  +2.85%  '

$ws.Range("D11").Value = '''59.87'
$ws.Range("E11").Value = '  +2.64%  '

$ws.Range("D12").Value = '''0.0939'
$ws.Range("E12").Value = '  -1.16%  '

$ws.Range("D13").Value = '''7.11'
$ws.Range("E13").Value = '  -1.31%  '

$ws.Range("E14").Value = '  -0.98%  '

$ws.Range("D15").Value = '2.554.51'
$ws.Range("E15").Value = '  +1.07%  '

$ws.Range("D16").Value = '''14.65'
$ws.Range("E16").Value = '  -0.82%  '

$ws.Range("D17").Value = '''0.872'
$ws.Range("E17").Value = '  -1.39%  '

$ws.Range("D18").Value = '2.220.19'
$ws.Range("E18").Value = '  +1.10%  '

$ws.Range("D19").Value = '41.866.37'
$ws.Range("E19").Value = '  +1.16%  '

$ws.Range("E20").Value = '  +1.07%  '

$ws.Range("D21").Value = '''6.22'
$ws.Range("E21").Value = '  -0.79%  '

$ws.Range("D22").Value = '''72.65'
$ws.Range("E22").Value = '  +0.91%  '

$ws.Range("D23").Value = '''232.07'
$ws.Range("E23").Value = '  -0.48%  '

$ws.Range("D24").Value = '''2.07'
$ws.Range("E24").Value = '  -1.03%  '

$ws.Range("D25").Value = '''3.93'
$ws.Range("E25").Value = '  +0.56%  '

$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '''11.37'
$ws.Range("E27").Value = '  -5.56%  '

$ws.Range("E28").Value = '  -4.31%  '

$ws.Range("E29").Value = '  -1.71%  '

$ws.Range("E30").Value = '  -1.97%  '

$ws.Range("D31").Value = '''167.08'
$ws.Range("E31").Value = '  -1.93%  '

$ws.Range("D32").Value = '''20.43'
$ws.Range("E32").Value = '  -1.26%  '

$ws.Range("D33").Value = '''0.0801'
$ws.Range("E33").Value = '  +9.46%  '

$ws.Range("D34").Value = '''5.93'
$ws.Range("E34").Value = '  +5.87%  '

$ws.Range("D35").Value = '''0.120'
$ws.Range("E35").Value = '  -1.11%  '

$ws.Range("E36").Value = '  -0.36%  '

$ws.Range("D37").Value = '''4.61'
$ws.Range("E37").Value = '  -0.52%  '

$ws.Range("D38").Value = '''4.11'
$ws.Range("E38").Value = '  +1.83%  '

$ws.Range("D39").Value = '''25.58'
$ws.Range("E39").Value = '  -3.23%  '

$ws.Range("E40").Value = '  +2.19%  '

$ws.Range("D41").Value = '''2.24'

$ws.Range("D42").Value = '''12.21'

$ws.Range("D43").Value = '''5.64'
$ws.Range("E43").Value = '  -2.85%  '

$ws.Range("D44").Value = '''5.07'
$ws.Range("E44").Value = '  +1.31%  '

$ws.Range("D45").Value = '''62.15'
$ws.Range("E45").Value = '  -3.11%  '

$ws.Range("D46").Value = '''0.200'
$ws.Range("E46").Value = '  -2.32%  '

$ws.Range("D47").Value = '''8.60'
$ws.Range("E47").Value = '  -1.42%  '

$ws.Range("E48").Value = '  -0.62%  '

$ws.Range("E49").Value = '  -0.42%  '

$ws.Range("E50").Value = '  +0.81%  '

$ws.Range("D51").Value = '''4.35'
$ws.Range("E51").Value = '  +1.04%  '
